$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix Species_name/Common_name/Category labels for the "Urophycis sp" / "unassigned" rows ---
# (rows were mislabeled; diff swaps the content between the two rows in each station block,
#  and capitalizes "unassigned" -> "Unassigned")
$ws.Range("A39").Value = "Unassigned"
$ws.Range("B39").Value = "Unassigned"
$ws.Range("C39").Value = "Unassigned"

$ws.Range("A40").Value = "Urophycis sp"
$ws.Range("B40").Value = "Red White or Spotted hake"
$ws.Range("C40").Value = "Teleost Fish"

$ws.Range("A78").Value = "Unassigned"
$ws.Range("B78").Value = "Unassigned"
$ws.Range("C78").Value = "Unassigned"

$ws.Range("A79").Value = "Urophycis sp"
$ws.Range("B79").Value = "Red White or Spotted hake"
$ws.Range("C79").Value = "Teleost Fish"

$ws.Range("A117").Value = "Unassigned"
$ws.Range("B117").Value = "Unassigned"
$ws.Range("C117").Value = "Unassigned"

$ws.Range("A118").Value = "Urophycis sp"
$ws.Range("B118").Value = "Red White or Spotted hake"
$ws.Range("C118").Value = "Teleost Fish"

# --- Recomputed relative-abundance (relab) values in column E ---
$ws.Range("E4").Value = [double]"0"
$ws.Range("E5").Value = [double]"0"
$ws.Range("E7").Value = [double]"0.02901290799401867"
$ws.Range("E9").Value = [double]"0"
$ws.Range("E11").Value = [double]"0.04035681751774634"
$ws.Range("E12").Value = [double]"0"
$ws.Range("E13").Value = [double]"0.005465701861432426"
$ws.Range("E18").Value = [double]"0"
$ws.Range("E19").Value = [double]"0.09575290902528316"
$ws.Range("E20").Value = [double]"0"
$ws.Range("E22").Value = [double]"0.01014076760454444"
$ws.Range("E23").Value = [double]"0"
$ws.Range("E24").Value = [double]"0"
$ws.Range("E25").Value = [double]"0.6262353689348756"
$ws.Range("E27").Value = [double]"0.001598459978343446"
$ws.Range("E28").Value = [double]"0.03330984341967309"
$ws.Range("E29").Value = [double]"0.002165655454529829"
$ws.Range("E31").Value = [double]"0.1302143311390316"
$ws.Range("E33").Value = [double]"0.001289080627696327"
$ws.Range("E34").Value = [double]"6.875096681047077E-05"
$ws.Range("E35").Value = [double]"0.002715663189013596"
$ws.Range("E38").Value = [double]"0"
$ws.Range("E39").Value = [double]"0"
$ws.Range("E40").Value = [double]"0.02167374228700091"
$ws.Range("E41").Value = [double]"0.08015624140397205"
$ws.Range("E43").Value = [double]"0"
$ws.Range("E44").Value = [double]"0"
$ws.Range("E47").Value = [double]"0.001283673506812639"
$ws.Range("E48").Value = [double]"0"
$ws.Range("E49").Value = [double]"0.003062478223395866"
$ws.Range("E51").Value = [double]"0"
$ws.Range("E53").Value = [double]"0.001760466523628762"
$ws.Range("E54").Value = [double]"0.001925510260218958"
$ws.Range("E55").Value = [double]"0.0120665309640388"
$ws.Range("E56").Value = [double]"0.0458821587720746"
$ws.Range("E57").Value = [double]"0"
$ws.Range("E58").Value = [double]"0.003704314976802186"
$ws.Range("E59").Value = [double]"0"
$ws.Range("E60").Value = [double]"0.0008618950688599145"
$ws.Range("E61").Value = [double]"0.1729474977535713"
$ws.Range("E62").Value = [double]"0"
$ws.Range("E63").Value = [double]"0"
$ws.Range("E64").Value = [double]"0.01452384881993728"
$ws.Range("E65").Value = [double]"0.0003117492802259266"
$ws.Range("E66").Value = [double]"0.02416973831398654"
$ws.Range("E67").Value = [double]"0.01454218701289175"
$ws.Range("E68").Value = [double]"0.03731822266233886"
$ws.Range("E69").Value = [double]"0.006583411270653389"
$ws.Range("E70").Value = [double]"0.02070381984559241"
$ws.Range("E71").Value = [double]"0.003172507381122664"
$ws.Range("E73").Value = [double]"0.01709119583356256"
$ws.Range("E74").Value = [double]"0.491188498285379"
$ws.Range("E75").Value = [double]"0.02099723093286388"
$ws.Range("E76").Value = [double]"0.001485393629311768"
$ws.Range("E77").Value = [double]"0"
$ws.Range("E78").Value = [double]"0.005464781500430948"
$ws.Range("E79").Value = [double]"0.01879664777832792"
$ws.Range("E81").Value = [double]"0.4732762888056007"
$ws.Range("E82").Value = [double]"0"
$ws.Range("E83").Value = [double]"0"
$ws.Range("E84").Value = [double]"0.02574075383636235"
$ws.Range("E87").Value = [double]"0"
$ws.Range("E90").Value = [double]"0"
$ws.Range("E92").Value = [double]"0.0411710628668411"
$ws.Range("E96").Value = [double]"0"
$ws.Range("E98").Value = [double]"0"
$ws.Range("E100").Value = [double]"0.2076939396082314"
$ws.Range("E101").Value = [double]"0"
$ws.Range("E102").Value = [double]"0"
$ws.Range("E106").Value = [double]"0.0004808712255144615"
$ws.Range("E107").Value = [double]"4.24298140159819E-05"
$ws.Range("E113").Value = [double]"0.2480588360087688"
$ws.Range("E116").Value = [double]"0"
$ws.Range("E117").Value = [double]"0.003380241849939891"
$ws.Range("E118").Value = [double]"0.000155575984725267"
